# Update the workbook's embedded build-version string everywhere it appears.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# ---- "About" sheet ----
$about = $wb.Worksheets.Item("About")

# A2: "Version: <version string>"
$a2 = $about.Range("A2")
$a2.Value = $a2.Value().Replace($oldVersion, $newVersion)

# A6: Recommended citation text containing the version string
$a6 = $about.Range("A6")
$a6.Value = $a6.Value().Replace($oldVersion, $newVersion)

# ---- "Boundaries and methane sources" sheet ----
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") holds the version string for each data row.
# Data rows run from row 2 through row 24.
for ($r = 2; $r -le 24; $r++) {
    $cell = $data.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value() -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
